# Auto-generated Excel COM-interop script to apply the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

$ws.Range('D2').Value = '29.209.38'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '1.866.60'
$ws.Range('E3').Value = '  -0.84%  '
Set-TextValue 'D4' '0.9999'
$ws.Range('E4').Value = '  -0.09%  '
Set-TextValue 'D5' '0.7108'
$ws.Range('E5').Value = '  -0.78%  '
Set-TextValue 'D6' '241.47'
$ws.Range('E6').Value = '  -0.34%  '
Set-TextValue 'D7' '1.000'
$ws.Range('E7').Value = '  -0.09%  '
Set-TextValue 'D8' '0.3115'
$ws.Range('E8').Value = '  -0.05%  '
Set-TextValue 'D9' '0.07663'
$ws.Range('E9').Value = '  -3.62%  '
Set-TextValue 'D10' '24.70'
$ws.Range('E10').Value = '  -2.52%  '
Set-TextValue 'D11' '0.08359'
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('D12').Value = '1.861.52'
$ws.Range('E12').Value = '  -0.01%  '
Set-TextValue 'D13' '5.220'
$ws.Range('E13').Value = '  -1.37%  '
Set-TextValue 'D14' '0.7112'
Set-TextValue 'D15' '91.30'
$ws.Range('D16').Value = '29.224.04'
$ws.Range('E16').Value = '  -1.08%  '
Set-TextValue 'D17' '5.943'
$ws.Range('E17').Value = '  -0.15%  '
Set-TextValue 'D18' '243.39'
$ws.Range('E18').Value = '  -1.24%  '
Set-TextValue 'D19' '0.000007817'
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').Value = '2.114.10'
$ws.Range('E20').Value = '  -0.44%  '
Set-TextValue 'D21' '13.10'
$ws.Range('E21').Value = '  -2.16%  '
Set-TextValue 'D22' '0.9990'
Set-TextValue 'D23' '7.856'
$ws.Range('E23').Value = '  -1.62%  '
Set-TextValue 'D24' '1.000'
$ws.Range('E24').Value = '  -0.08%  '
Set-TextValue 'D25' '0.1594'
$ws.Range('E25').Value = '  -1.55%  '
Set-TextValue 'D26' '163.10'
$ws.Range('E26').Value = '  -0.19%  '
Set-TextValue 'D27' '8.942'
Set-TextValue 'D28' '18.47'
$ws.Range('E28').Value = '  +0.55%  '
Set-TextValue 'D29' '1.500'
$ws.Range('E29').Value = '  -0.01%  '
Set-TextValue 'D30' '1.316'
$ws.Range('E30').Value = '  -3.15%  '
Set-TextValue 'D31' '4.398'
$ws.Range('E31').Value = '  -0.05%  '
Set-TextValue 'D32' '4.247'
$ws.Range('E32').Value = '  +3.18%  '
Set-TextValue 'D33' '0.05159'
$ws.Range('E33').Value = '  -2.80%  '
Set-TextValue 'D34' '0.7965'
$ws.Range('E34').Value = '  +9.35%  '
Set-TextValue 'D35' '1.912'
$ws.Range('E35').Value = '  -2.70%  '
Set-TextValue 'D37' '2.685'
$ws.Range('E37').Value = '  +0.17%  '
Set-TextValue 'D38' '0.01850'
$ws.Range('E38').Value = '  -1.19%  '
Set-TextValue 'D39' '2.708'
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').Value = '1.159.37'
$ws.Range('E40').Value = '  -6.00%  '
Set-TextValue 'D41' '6.300'
$ws.Range('E41').Value = '  +1.25%  '
Set-TextValue 'D42' '0.8952'
$ws.Range('E42').Value = '  -1.81%  '
Set-TextValue 'D43' '73.10'
$ws.Range('E43').Value = '  -1.16%  '
Set-TextValue 'D44' '0.9994'
$ws.Range('E44').Value = '  -0.19%  '
Set-TextValue 'D45' '103.01'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('D46').Value = '2.011.48'
$ws.Range('E46').Value = '  -0.86%  '
Set-TextValue 'D47' '0.5187'
$ws.Range('E47').Value = '  -1.84%  '
Set-TextValue 'D48' '1.779'
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D49' '0.00000000120'
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '9.336'
$ws.Range('E50').Value = '  -0.05%  '
Set-TextValue 'D51' '0.4291'
$ws.Range('E51').Value = '  -1.17%  '
